$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effort")

# New test case tc10 entries: system load estimation
# Shared-string pool order follows first-write order, so populate the
# text cells in the same order as the original author: "System load
# estimation" (D48) first, then "Minor changes ..." (D47), then the
# tc10 note (D49).
$ws.Range("D48").Value = "System load estimation"
$ws.Range("D47").Value = "Minor changes on documentation and setup. System load estimation"
$ws.Range("D49").Value = "tc: System load estimation put to operation, validated by test case tc10"

$ws.Range("A47").Value = (Get-Date -Year 2012 -Month 11 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B47").Value = 2

$ws.Range("A48").Value = (Get-Date -Year 2012 -Month 11 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B48").Value = 2

$ws.Range("A49").Value = (Get-Date -Year 2012 -Month 11 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B49").Value = 2

$ws.Range("A50").Select()
